$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: only C changes
$ws.Range("C2").Value = 0.6264

# Row 3: B and C change
$ws.Range("B3").Value = 0.1465
$ws.Range("C3").Value = 0.454

# Rows 4-24: only C changes
$ws.Range("C4").Value = 0.2719
$ws.Range("C5").Value = 0.5312
$ws.Range("C6").Value = 0.7237
$ws.Range("C7").Value = 0.7379
$ws.Range("C8").Value = 0.3672
$ws.Range("C9").Value = 0.1072
$ws.Range("C10").Value = -0.0074
$ws.Range("C11").Value = -0.3166
$ws.Range("C12").Value = -0.4143
$ws.Range("C13").Value = 0.4814
$ws.Range("C14").Value = 0.345
$ws.Range("C15").Value = 0.4854
$ws.Range("C16").Value = 1.4614
$ws.Range("C17").Value = 1.3536
$ws.Range("C18").Value = 1.8308
$ws.Range("C19").Value = 1.7431
$ws.Range("C20").Value = 1.1392
$ws.Range("C21").Value = 1.008
$ws.Range("C22").Value = 0.6159
$ws.Range("C23").Value = 0.6297
$ws.Range("C24").Value = 0.4218
